# "carga del excel lista" - update the afiliados/personas data on Hoja1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 2: Argenis Garcia
$ws.Range("A2").Value = "V21443181"
$ws.Range("B2").Value = "Argenis"
$ws.Range("C2").Value = "Garcia"
$ws.Range("D2").Value = "04/02/1993"
$ws.Range("E2").Value = "M"
$ws.Range("F2").Value = "4127849712 454548"

# Row 3: Ivonne Ortega
$ws.Range("A3").Value = "V20699411"
$ws.Range("B3").Value = "Ivonne"
$ws.Range("C3").Value = "Ortega"
$ws.Range("D3").Value = "10/28/1993"
$ws.Range("E3").Value = "F"
$ws.Range("F3").Value = 123456

# Row 4: Priscilla BlaBla
$ws.Range("A4").Value = "V123456789"
$ws.Range("B4").Value = "Priscilla"
$ws.Range("C4").Value = "BlaBla"
$ws.Range("D4").Value = "02/12/1990"
$ws.Range("E4").Value = "F"
$ws.Range("F4").Value = 23435234
